# Move the identified "bad zombie" sapling.id values from their current
# rows to the bottom of the list (in their original relative order),
# shifting everything else up to fill the gaps.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# sapling.id values in data/zombie8.xlsx previously flagged as bad zombies
$zombieValues = @(9, 97, 142, 1314, 1331)

$firstDataRow = 2
$lastDataRow = 101

# Track how many of each zombie value still need to be pulled out, in case
# of duplicates (none expected here, but keep this robust).
$zombieCounts = @{}
foreach ($z in $zombieValues) {
    if ($zombieCounts.ContainsKey($z)) {
        $zombieCounts[$z] = $zombieCounts[$z] + 1
    } else {
        $zombieCounts[$z] = 1
    }
}

# Read the current column of sapling.id values (row 2..101) into an array.
$values = @()
for ($r = $firstDataRow; $r -le $lastDataRow; $r++) {
    $cellValue = [int]$ws.Cells.Item($r, 1).Value2
    $values += $cellValue
}

# Pull the zombie values out, preserving the relative order of everything
# else, then append the zombies themselves at the end in the order they
# originally appeared.
$remaining = @()
$pulled = @()

foreach ($v in $values) {
    if ($zombieCounts.ContainsKey($v) -and $zombieCounts[$v] -gt 0) {
        $zombieCounts[$v] = $zombieCounts[$v] - 1
        $pulled += $v
    } else {
        $remaining += $v
    }
}

$newOrder = @()
$newOrder += $remaining
$newOrder += $pulled

# Write the reordered values back into A2:A101.
for ($i = 0; $i -lt $newOrder.Count; $i++) {
    $row = $firstDataRow + $i
    $ws.Cells.Item($row, 1).Value = $newOrder[$i]
}
